$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("SpecialPrices")
$ws3.Columns.Item(4).ColumnWidth = 194.16666666666666
